$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Rows 2-28: price/volume updates only (coin name & link unchanged)
Set-TextCell "D2" "69.735.67"
Set-TextCell "E2" "  +2.18%  "

Set-TextCell "D3" "3.939.52"
Set-TextCell "E3" "  +0.66%  "

Set-TextCell "E4" "  +0.06%  "

Set-TextCell "D5" "532.38"
Set-TextCell "E5" "  +9.02%  "

Set-TextCell "D6" "146.51"
Set-TextCell "E6" "  -0.32%  "

Set-TextCell "E7" "  +0.02%  "

Set-TextCell "E8" "  -0.04%  "

Set-TextCell "D9" "0.733"
Set-TextCell "E9" "  +0.63%  "

Set-TextCell "E10" "  +4.89%  "

Set-TextCell "D11" "0.0000343"
Set-TextCell "E11" "  -0.57%  "

Set-TextCell "D12" "42.98"
Set-TextCell "E12" "  +0.08%  "

Set-TextCell "D13" "10.49"
Set-TextCell "E13" "  -2.26%  "

Set-TextCell "D14" "4.568.63"
Set-TextCell "E14" "  +0.80%  "

Set-TextCell "D15" "3.933.45"
Set-TextCell "E15" "  +0.91%  "

Set-TextCell "E16" "  -0.06%  "

Set-TextCell "E17" "  -0.13%  "

Set-TextCell "E18" "  +6.99%  "

Set-TextCell "D19" "19.94"
Set-TextCell "E19" "  +0.51%  "

Set-TextCell "D20" "69.535.90"
Set-TextCell "E20" "  +1.81%  "

Set-TextCell "D21" "435.45"
Set-TextCell "E21" "  +0.77%  "

Set-TextCell "D22" "3.41"
Set-TextCell "E22" "  -4.39%  "

Set-TextCell "E23" "  -2.32%  "

Set-TextCell "D24" "88.60"
Set-TextCell "E24" "  +1.12%  "

Set-TextCell "D25" "4.12"
Set-TextCell "E25" "  +14.19%  "

Set-TextCell "D26" "11.97"
Set-TextCell "E26" "  +3.88%  "

Set-TextCell "D27" "11.01"
Set-TextCell "E27" "  -2.60%  "

Set-TextCell "D28" "36.74"
Set-TextCell "E28" "  -3.72%  "

# Rows 29-51: the LEO coin (previously row 29) drops out of the list, and
# every subsequent coin shifts up by one row; a new coin (Monero) is
# appended at the bottom (row 51). Column A (rank index) is untouched.

$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D29" "711.71"
Set-TextCell "E29" "  -2.02%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D30" "13.40"
Set-TextCell "E30" "  -2.55%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.128"
Set-TextCell "E31" "  -1.61%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D32" "2.87"
Set-TextCell "E32" "  -1.59%  "

$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D33" "68.67"
Set-TextCell "E33" "  +13.93%  "

$ws.Range("B34").Value = "TheGraph"
$ws.Range("C34").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D34" "0.447"
Set-TextCell "E34" "  +9.23%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D35" "6.08"
Set-TextCell "E35" "  -3.23%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D36" "0.0₃0872"
Set-TextCell "E36" "  +0.36%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D37" "40.69"
Set-TextCell "E37" "  -2.62%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D38" "0.150"
Set-TextCell "E38" "  +1.08%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D39" "0.999"
Set-TextCell "E39" "  -0.06%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D40" "1.00"
Set-TextCell "E40" "  +0.09%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D41" "0.0485"
Set-TextCell "E41" "  +0.88%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D42" "2.83"
Set-TextCell "E42" "  -4.78%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D43" "3.10"
Set-TextCell "E43" "  +6.19%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell "D44" "3.02"
Set-TextCell "E44" "  -4.57%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D45" "3.27"
Set-TextCell "E45" "  +16.53%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D46" "0.143"
Set-TextCell "E46" "  +1.47%  "

# Row 47 (ApeXProtocol) keeps its name/link, only price/volume updates
Set-TextCell "D47" "3.39"
Set-TextCell "E47" "  +2.60%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.0₆0364"
Set-TextCell "E48" "  +5.34%  "

$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D49" "3.36"
Set-TextCell "E49" "  -1.53%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D50" "2.10"
Set-TextCell "E50" "  -1.17%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D51" "144.89"
Set-TextCell "E51" "  +0.12%  "
